$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells P1, Q1 with values, copying style/format from O1 (the last existing header cell)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update I, K, M, O columns for rows 2-25 (swap values between the 1s and 2s),
# and populate new P, Q columns with value 2 for each row
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I column
    $ws.Cells.Item($r, 11).Value = 1  # K column
    $ws.Cells.Item($r, 13).Value = 2  # M column
    $ws.Cells.Item($r, 15).Value = 1  # O column
    $ws.Cells.Item($r, 16).Value = 2  # P column
    $ws.Cells.Item($r, 17).Value = 2  # Q column
}
